# Insert a new data row at row 690 (pushes existing rows 690-783 down to 691-784)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(690).Insert()

$ws.Range("A690").Value = 10
$ws.Range("B690").Value = "Vega Modelo de Temuco"
$ws.Range("C690").Value = "La Araucanía"
$ws.Range("D690").Value = 44918
$ws.Range("E690").Value = 9
$ws.Range("F690").Value = "Fruta"
$ws.Range("G690").Value = 100102
$ws.Range("H690").Value = "Cítricos"
$ws.Range("I690").Value = 100102004
$ws.Range("J690").Value = "Mandarina"
$ws.Range("K690").Value = "Murcott"
$ws.Range("L690").Value = "Tercera"
$ws.Range("M690").Value = 2
$ws.Range("N690").Value = 160000
$ws.Range("O690").Value = 160000
$ws.Range("P690").Value = 160000
$ws.Range("Q690").Value = "$/bins (450 kilos)"
$ws.Range("R690").Value = "Región de O'Higgins"
$ws.Range("S690").Value = 356
$ws.Range("T690").Value = 450
